# StructureDefinition-benefit-status.xlsx update
# - Bump Version from 5.0.0 to 6.0.0
# - Bump Date to the new publication timestamp
# - Fill in Publisher ("Alvearie Team")
# - Replace the stray duplicate "Contact" row with a "Jurisdiction" row,
#   and drop the now-redundant extra "Contact" row entirely (net -1 row)
# - Mirror the profile-specific Short/Definition text onto the root
#   Extension row of the element-detail sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 ("Metadata") ---------------------------------------------

# Version
$ws1.Range("B3").Value = "6.0.0"

# Date
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a second "Contact" row; turn it into "Jurisdiction"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was the original duplicate "Contact" / "No display for ContactDetail"
# row - remove it outright so everything below shifts up by one row.
$ws1.Rows.Item(11).Delete()

# --- Sheet 2 (element detail grid) ------------------------------------

# The root "Extension" row's Short/Definition columns (K/L) should carry
# the profile's own Title/Description instead of the generic placeholder.
$ws2.Range("K2").Value = "Benefit Status"
$ws2.Range("L2").Value = "HIPAA standard code for the benefit status"
